$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110; this shifts the existing rows
# 110..248 down to 111..249 (dimension grows from R248 to R249).
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly record.
# The "constant" columns (A,B,C,E,F,G,H,I,O,R) are identical for every
# data row in this sheet.
$ws.Range("A110").Value = 8
$ws.Range("B110").Value = "Terminal La Palmera de La Serena"
$ws.Range("C110").Value = "Coquimbo"
$ws.Range("D110").Value = 44848
$ws.Range("E110").Value = 4
$ws.Range("F110").Value = 100112037
$ws.Range("G110").Value = "Cebollín"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 1400
$ws.Range("L110").Value = 1600
$ws.Range("M110").Value = 1500
$ws.Range("N110").Value = "$/paquete 6 unidades"
$ws.Range("O110").Value = "Provincia del Elquí"
$ws.Range("P110").Value = 250
$ws.Range("Q110").Value = 6
$ws.Range("R110").Value = "Hortaliza"
